$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 208, shifting existing rows 208-227 down to 209-228
$ws.Rows.Item(208).Insert()

# The newly inserted row 208 is currently blank; copy formatting from row 209 (previously row 208)
# for consistency, then set its values explicitly.

$row = 208
$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44826
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(209, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = "Arándano (blue)"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 50
$ws.Cells.Item($row, 14).Value = 12000
$ws.Cells.Item($row, 15).Value = 12000
$ws.Cells.Item($row, 16).Value = 12000
$ws.Cells.Item($row, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 8000
$ws.Cells.Item($row, 20).Value = 1.5

$wb.Save()
